# Apply edits to the amref workbook as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header (row 1) text renames -----------------------------------------
$ws.Range("C1").Value  = "GDP"
$ws.Range("E1").Value  = "Budget_Previous_Year"
$ws.Range("F1").Value  = "LatinAmerica"
$ws.Range("G1").Value  = "Africa"
$ws.Range("H1").Value  = "Confessional"
$ws.Range("I1").Value  = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C (GDP) data updates, rows 2-31 -------------------------------
$ws.Range("C2").Value  = 743.403784726004
$ws.Range("C3").Value  = 341.5541149051794
$ws.Range("C4").Value  = 951.6879611168786
$ws.Range("C5").Value  = 1401.47747416771
$ws.Range("C6").Value  = 815.8736791314819
$ws.Range("C7").Value  = 777.227218443918
$ws.Range("C8").Value  = 369.2024078290272
$ws.Range("C9").Value  = 982.980837581714
$ws.Range("C10").Value = 1591.56825353313
$ws.Range("C11").Value = 864.5379000312432
$ws.Range("C12").Value = 389.9389667216314
$ws.Range("C13").Value = 1000.829216794104
$ws.Range("C14").Value = 788.439151581443
$ws.Range("C15").Value = 1745.10167474004
$ws.Range("C16").Value = 1291.622214254295
$ws.Range("C17").Value = 419.1838602515346
$ws.Range("C18").Value = 1032.277326842402
$ws.Range("C19").Value = 817.1226340535979
$ws.Range("C20").Value = 1778.60982580794
$ws.Range("C21").Value = 1291.415042301529
$ws.Range("C22").Value = 449.4203771491282
$ws.Range("C23").Value = 741.0381351906716
$ws.Range("C24").Value = 1060.095015975378
$ws.Range("C25").Value = 482.6390663355013
$ws.Range("C26").Value = 750.4706590411453
$ws.Range("C27").Value = 1093.134170274031
$ws.Range("C28").Value = 514.0573067519859
$ws.Range("C29").Value = 809.9545825255682
$ws.Range("C30").Value = 1129.713195979213
$ws.Range("C31").Value = 1431.756130822538
